$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the existing row 493, shifting the
# old rows 494-579 down to 496-581 (dimension grows from A1:R579 to A1:R581).
$ws.Rows("494:495").Insert()

# Populate the two newly inserted rows with the new weekly price records.
$ws.Range("A494").Value = 9
$ws.Range("B494").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C494").Value = "Metropolitana"
$ws.Range("D494").Value = 44637
$ws.Range("E494").Value = 13
$ws.Range("F494").Value = 100112008
$ws.Range("G494").Value = "Coliflor"
$ws.Range("H494").Value = "Sin especificar"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 2500
$ws.Range("K494").Value = 900
$ws.Range("L494").Value = 1000
$ws.Range("M494").Value = 950
$ws.Range("N494").Value = "$/unidad"
$ws.Range("O494").Value = "Región Metropolitana"
$ws.Range("P494").Value = 950
$ws.Range("Q494").Value = 1
$ws.Range("R494").Value = "Hortaliza"

$ws.Range("A495").Value = 9
$ws.Range("B495").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C495").Value = "Metropolitana"
$ws.Range("D495").Value = 44637
$ws.Range("E495").Value = 13
$ws.Range("F495").Value = 100112008
$ws.Range("G495").Value = "Coliflor"
$ws.Range("H495").Value = "Sin especificar"
$ws.Range("I495").Value = "Segunda"
$ws.Range("J495").Value = 970
$ws.Range("K495").Value = 800
$ws.Range("L495").Value = 800
$ws.Range("M495").Value = 800
$ws.Range("N495").Value = "$/unidad"
$ws.Range("O495").Value = "Región Metropolitana"
$ws.Range("P495").Value = 800
$ws.Range("Q495").Value = 1
$ws.Range("R495").Value = "Hortaliza"
